# Added New Mac-Address and Document Types
# Appends five new detail rows (157-161) to the
# master-reg_center_machine_device_h worksheet, mirroring the pattern of the
# existing rows (regcntr_id/machine_id pairs with incrementing device_id),
# then updates the sheet's selection and switches the workbook to manual
# calculation mode.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# lang_code="eng", is_active=TRUE, cr_by="superadmin", cr_dtimes/eff_dtimes="now()"
$regcntrId = 10002
$machineId = 10032
$deviceIds = @(3000176, 3000177, 3000178, 3000179, 3000180)

$row = 157
foreach ($deviceId in $deviceIds) {
    $ws.Cells.Item($row, 1).Value = $regcntrId
    $ws.Cells.Item($row, 2).Value = $machineId
    $ws.Cells.Item($row, 3).Value = $deviceId
    $ws.Cells.Item($row, 4).Value = "eng"
    $ws.Cells.Item($row, 5).Value = $true
    $ws.Cells.Item($row, 6).Value = "superadmin"
    $ws.Cells.Item($row, 7).Value = "now()"
    $ws.Cells.Item($row, 8).Value = "now()"
    $row = $row + 1
}

# Move the visible selection down to the newly added data, like the author
# scrolled to and selected the first new row's lang_code cell.
$ws.Range("D157").Select()

# Switch calculation mode to manual (<calcPr calcMode="manual"/>).
$excel.Calculation = -4135
